$wb = $excel.ActiveWorkbook

# --- Sheet "Overview" ---
$ws = $wb.Worksheets.Item("Overview")
$ws.Rows.Item(9).Insert()
$ws.Range("A9").Value = "ef12ab35-0a62-4357-af99-1bc873e0fcc3.md"
$ws.Range("B9").Value = "Ready for handoff"
$ws.Range("C9").Value = "Ready for handoff"

# Rebuild hyperlinks (engine only supports whole-sheet clear)
$ws.Range("A1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/bdc60e2bda5902d2ba3712a83dc583abec29b99c/e2e/f09010b7-c32b-477d-9e73-4999517e5677.md", "", "", "f09010b7-c32b-477d-9e73-4999517e5677.md")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/74dc4bd4941e3ccea240f7cdac99e11dd5b9d2b2/e2e/11c04f05-8be1-4b0a-8bcc-3413fe5e1340.md", "", "", "11c04f05-8be1-4b0a-8bcc-3413fe5e1340.md")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/74dc4bd4941e3ccea240f7cdac99e11dd5b9d2b2/e2e/2f541edf-1173-465c-a5ba-f1619df9c157.md", "", "", "2f541edf-1173-465c-a5ba-f1619df9c157.md")
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/0981438d1a06861daaf8a92ebc4c18b36fc4f35a/e2e/4ea1af5a-5d3c-42fa-8e63-933cff256c94.md", "", "", "4ea1af5a-5d3c-42fa-8e63-933cff256c94.md")
$ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/8290aabd96d7fe3d5f38f12e392d1c2a4d307999/e2e/85077eb4-7547-4b81-859a-c9cfe3701a6d.md", "", "", "85077eb4-7547-4b81-859a-c9cfe3701a6d.md")
$ws.Hyperlinks.Add($ws.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/15efe549cc8227c412e1572d8944ebe5eb84aabe/e2e/82429938-e6c9-4c64-aeed-848f6d261f77.md", "", "", "82429938-e6c9-4c64-aeed-848f6d261f77.md")
$ws.Hyperlinks.Add($ws.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/0e585a68cdc383e277f5da1c49d03c82fdcdc23c/e2e/978a0a30-878a-4dbc-a5ac-038765e6fc83.md", "", "", "978a0a30-878a-4dbc-a5ac-038765e6fc83.md")
$ws.Hyperlinks.Add($ws.Range("A9"), "https://github.com/OpenLocalizationTest/oltest/blob/bdc60e2bda5902d2ba3712a83dc583abec29b99c/e2e/ef12ab35-0a62-4357-af99-1bc873e0fcc3.md", "", "", "ef12ab35-0a62-4357-af99-1bc873e0fcc3.md")
$ws.Hyperlinks.Add($ws.Range("A10"), "https://github.com/OpenLocalizationTest/oltest/blob/bdc60e2bda5902d2ba3712a83dc583abec29b99c/.localization-config", "", "", ".localization-config")

# --- Sheet "zh-cn" ---
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Rows.Item(9).Insert()
$ws.Range("A9").Value = "ef12ab35-0a62-4357-af99-1bc873e0fcc3.md"
$ws.Range("B9").Value = "Ready for handoff"
$ws.Range("C9").Value = "ef12ab35-0a62-4357-af99-1bc873e0fcc3.b0a908d97a0c162733732f1a2a4d7c80cbd57c92.zh-cn.xlf"
$ws.Range("D9").Value = "2016-02-25 03:48:11"
$ws.Range("G9").Value = "0001-01-01 00:00:00"
$ws.Range("H9").Value = "Include"

# Rebuild hyperlinks (engine only supports whole-sheet clear)
$ws.Range("A1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/bdc60e2bda5902d2ba3712a83dc583abec29b99c/e2e/f09010b7-c32b-477d-9e73-4999517e5677.md", "", "", "f09010b7-c32b-477d-9e73-4999517e5677.md")
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8e7fdf120d370e450aa27db1530b0e972d1dbf87/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/f09010b7-c32b-477d-9e73-4999517e5677.dd782c6cf69313603a0c6c94730ebba3604cb6dd.zh-cn.xlf", "", "", "f09010b7-c32b-477d-9e73-4999517e5677.dd782c6cf69313603a0c6c94730ebba3604cb6dd.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/a7a4c2455b58dd982d79dec2fd64556083904819/e2e/f09010b7-c32b-477d-9e73-4999517e5677.md", "", "", "f09010b7-c32b-477d-9e73-4999517e5677.md")
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/b307e8d69eb953f385417d0e73be7583bd3448e0/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/f09010b7-c32b-477d-9e73-4999517e5677.dd782c6cf69313603a0c6c94730ebba3604cb6dd.zh-cn.xlf", "", "", "f09010b7-c32b-477d-9e73-4999517e5677.dd782c6cf69313603a0c6c94730ebba3604cb6dd.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/74dc4bd4941e3ccea240f7cdac99e11dd5b9d2b2/e2e/11c04f05-8be1-4b0a-8bcc-3413fe5e1340.md", "", "", "11c04f05-8be1-4b0a-8bcc-3413fe5e1340.md")
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/64ea5a71aaec06356e08edd3c58c5678b97d1f1c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/11c04f05-8be1-4b0a-8bcc-3413fe5e1340.a2317b3789fa00c725c571ba2d40b45bbf3575c6.zh-cn.xlf", "", "", "11c04f05-8be1-4b0a-8bcc-3413fe5e1340.a2317b3789fa00c725c571ba2d40b45bbf3575c6.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/74dc4bd4941e3ccea240f7cdac99e11dd5b9d2b2/e2e/2f541edf-1173-465c-a5ba-f1619df9c157.md", "", "", "2f541edf-1173-465c-a5ba-f1619df9c157.md")
$ws.Hyperlinks.Add($ws.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/64ea5a71aaec06356e08edd3c58c5678b97d1f1c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/2f541edf-1173-465c-a5ba-f1619df9c157.17f91ae2f4592db3b8c1784bd841d03852931f23.zh-cn.xlf", "", "", "2f541edf-1173-465c-a5ba-f1619df9c157.17f91ae2f4592db3b8c1784bd841d03852931f23.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/0981438d1a06861daaf8a92ebc4c18b36fc4f35a/e2e/4ea1af5a-5d3c-42fa-8e63-933cff256c94.md", "", "", "4ea1af5a-5d3c-42fa-8e63-933cff256c94.md")
$ws.Hyperlinks.Add($ws.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dae417aac267c6e53244d46e59cdb905672969df/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/4ea1af5a-5d3c-42fa-8e63-933cff256c94.fb647b1cf0a49b43fb46f22842a039fc7dba17f1.zh-cn.xlf", "", "", "4ea1af5a-5d3c-42fa-8e63-933cff256c94.fb647b1cf0a49b43fb46f22842a039fc7dba17f1.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/8290aabd96d7fe3d5f38f12e392d1c2a4d307999/e2e/85077eb4-7547-4b81-859a-c9cfe3701a6d.md", "", "", "85077eb4-7547-4b81-859a-c9cfe3701a6d.md")
$ws.Hyperlinks.Add($ws.Range("C6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c245c04ea8208cba0a07125583538087287d502d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/85077eb4-7547-4b81-859a-c9cfe3701a6d.95d5e296c59cf65950886e227155318594533518.zh-cn.xlf", "", "", "85077eb4-7547-4b81-859a-c9cfe3701a6d.95d5e296c59cf65950886e227155318594533518.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("E6"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/df53c81c9edfa118e900d173011e0e88e3dfbc2a/e2e/85077eb4-7547-4b81-859a-c9cfe3701a6d.md", "", "", "85077eb4-7547-4b81-859a-c9cfe3701a6d.md")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/5806f51c6301d3be3f0b95381f2bb8d829e1b58e/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/85077eb4-7547-4b81-859a-c9cfe3701a6d.95d5e296c59cf65950886e227155318594533518.zh-cn.xlf", "", "", "85077eb4-7547-4b81-859a-c9cfe3701a6d.95d5e296c59cf65950886e227155318594533518.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/15efe549cc8227c412e1572d8944ebe5eb84aabe/e2e/82429938-e6c9-4c64-aeed-848f6d261f77.md", "", "", "82429938-e6c9-4c64-aeed-848f6d261f77.md")
$ws.Hyperlinks.Add($ws.Range("C7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2310fbc175e5ad65f4f6fffe72b61a9af9ec9fee/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/82429938-e6c9-4c64-aeed-848f6d261f77.5ec01e0eb01a8d1518aac1cac8238bee3e6f1337.zh-cn.xlf", "", "", "82429938-e6c9-4c64-aeed-848f6d261f77.5ec01e0eb01a8d1518aac1cac8238bee3e6f1337.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/0e585a68cdc383e277f5da1c49d03c82fdcdc23c/e2e/978a0a30-878a-4dbc-a5ac-038765e6fc83.md", "", "", "978a0a30-878a-4dbc-a5ac-038765e6fc83.md")
$ws.Hyperlinks.Add($ws.Range("C8"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ad05b521ecebd8b5bae80d61e9d81203bbb388eb/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/978a0a30-878a-4dbc-a5ac-038765e6fc83.c933c29170ae9c133805761f650ea2014570545e.zh-cn.xlf", "", "", "978a0a30-878a-4dbc-a5ac-038765e6fc83.c933c29170ae9c133805761f650ea2014570545e.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A9"), "https://github.com/OpenLocalizationTest/oltest/blob/bdc60e2bda5902d2ba3712a83dc583abec29b99c/e2e/ef12ab35-0a62-4357-af99-1bc873e0fcc3.md", "", "", "ef12ab35-0a62-4357-af99-1bc873e0fcc3.md")
$ws.Hyperlinks.Add($ws.Range("C9"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a1f4c6e8b23d059172a6e4c8d0b3f7a9510cde6b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/ef12ab35-0a62-4357-af99-1bc873e0fcc3.b0a908d97a0c162733732f1a2a4d7c80cbd57c92.zh-cn.xlf", "", "", "ef12ab35-0a62-4357-af99-1bc873e0fcc3.b0a908d97a0c162733732f1a2a4d7c80cbd57c92.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A10"), "https://github.com/OpenLocalizationTest/oltest/blob/bdc60e2bda5902d2ba3712a83dc583abec29b99c/.localization-config", "", "", ".localization-config")

# --- Sheet "de-de" ---
$ws = $wb.Worksheets.Item("de-de")
$ws.Rows.Item(9).Insert()
$ws.Range("A9").Value = "ef12ab35-0a62-4357-af99-1bc873e0fcc3.md"
$ws.Range("B9").Value = "Ready for handoff"
$ws.Range("C9").Value = "ef12ab35-0a62-4357-af99-1bc873e0fcc3.b0a908d97a0c162733732f1a2a4d7c80cbd57c92.de-de.xlf"
$ws.Range("D9").Value = "2016-02-25 03:48:23"
$ws.Range("G9").Value = "0001-01-01 00:00:00"
$ws.Range("H9").Value = "Include"

# Rebuild hyperlinks (engine only supports whole-sheet clear)
$ws.Range("A1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/bdc60e2bda5902d2ba3712a83dc583abec29b99c/e2e/f09010b7-c32b-477d-9e73-4999517e5677.md", "", "", "f09010b7-c32b-477d-9e73-4999517e5677.md")
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0959cfa18cec7ead331c5e98499021dc18b08b4c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/f09010b7-c32b-477d-9e73-4999517e5677.dd782c6cf69313603a0c6c94730ebba3604cb6dd.de-de.xlf", "", "", "f09010b7-c32b-477d-9e73-4999517e5677.dd782c6cf69313603a0c6c94730ebba3604cb6dd.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/4ca943976a1fa2eac246180ebbfc76d9b3db9fb1/e2e/f09010b7-c32b-477d-9e73-4999517e5677.md", "", "", "f09010b7-c32b-477d-9e73-4999517e5677.md")
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a8791d5d1e0b9140d7baff2b9865544fe80b06a0/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/f09010b7-c32b-477d-9e73-4999517e5677.dd782c6cf69313603a0c6c94730ebba3604cb6dd.de-de.xlf", "", "", "f09010b7-c32b-477d-9e73-4999517e5677.dd782c6cf69313603a0c6c94730ebba3604cb6dd.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/74dc4bd4941e3ccea240f7cdac99e11dd5b9d2b2/e2e/11c04f05-8be1-4b0a-8bcc-3413fe5e1340.md", "", "", "11c04f05-8be1-4b0a-8bcc-3413fe5e1340.md")
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5879610df7cbe5c9a3cf2c5a5de93c72947ebe24/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/11c04f05-8be1-4b0a-8bcc-3413fe5e1340.a2317b3789fa00c725c571ba2d40b45bbf3575c6.de-de.xlf", "", "", "11c04f05-8be1-4b0a-8bcc-3413fe5e1340.a2317b3789fa00c725c571ba2d40b45bbf3575c6.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/74dc4bd4941e3ccea240f7cdac99e11dd5b9d2b2/e2e/2f541edf-1173-465c-a5ba-f1619df9c157.md", "", "", "2f541edf-1173-465c-a5ba-f1619df9c157.md")
$ws.Hyperlinks.Add($ws.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5879610df7cbe5c9a3cf2c5a5de93c72947ebe24/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/2f541edf-1173-465c-a5ba-f1619df9c157.17f91ae2f4592db3b8c1784bd841d03852931f23.de-de.xlf", "", "", "2f541edf-1173-465c-a5ba-f1619df9c157.17f91ae2f4592db3b8c1784bd841d03852931f23.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/0981438d1a06861daaf8a92ebc4c18b36fc4f35a/e2e/4ea1af5a-5d3c-42fa-8e63-933cff256c94.md", "", "", "4ea1af5a-5d3c-42fa-8e63-933cff256c94.md")
$ws.Hyperlinks.Add($ws.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d33491f4064a89ef7721a085337dc33329782522/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/4ea1af5a-5d3c-42fa-8e63-933cff256c94.fb647b1cf0a49b43fb46f22842a039fc7dba17f1.de-de.xlf", "", "", "4ea1af5a-5d3c-42fa-8e63-933cff256c94.fb647b1cf0a49b43fb46f22842a039fc7dba17f1.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/8290aabd96d7fe3d5f38f12e392d1c2a4d307999/e2e/85077eb4-7547-4b81-859a-c9cfe3701a6d.md", "", "", "85077eb4-7547-4b81-859a-c9cfe3701a6d.md")
$ws.Hyperlinks.Add($ws.Range("C6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7b53aeeea3bbce3b441436a62d5dd158747b32e7/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/85077eb4-7547-4b81-859a-c9cfe3701a6d.95d5e296c59cf65950886e227155318594533518.de-de.xlf", "", "", "85077eb4-7547-4b81-859a-c9cfe3701a6d.95d5e296c59cf65950886e227155318594533518.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("E6"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/046d2e953a54f82b19346662f486f83da6061f4b/e2e/85077eb4-7547-4b81-859a-c9cfe3701a6d.md", "", "", "85077eb4-7547-4b81-859a-c9cfe3701a6d.md")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/ac1b13dc57e2c4f92fc0daa77c537df483d77594/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/85077eb4-7547-4b81-859a-c9cfe3701a6d.95d5e296c59cf65950886e227155318594533518.de-de.xlf", "", "", "85077eb4-7547-4b81-859a-c9cfe3701a6d.95d5e296c59cf65950886e227155318594533518.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/15efe549cc8227c412e1572d8944ebe5eb84aabe/e2e/82429938-e6c9-4c64-aeed-848f6d261f77.md", "", "", "82429938-e6c9-4c64-aeed-848f6d261f77.md")
$ws.Hyperlinks.Add($ws.Range("C7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2cb5500c1682bb75512eb318fa0d8b55276ae957/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/82429938-e6c9-4c64-aeed-848f6d261f77.5ec01e0eb01a8d1518aac1cac8238bee3e6f1337.de-de.xlf", "", "", "82429938-e6c9-4c64-aeed-848f6d261f77.5ec01e0eb01a8d1518aac1cac8238bee3e6f1337.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/0e585a68cdc383e277f5da1c49d03c82fdcdc23c/e2e/978a0a30-878a-4dbc-a5ac-038765e6fc83.md", "", "", "978a0a30-878a-4dbc-a5ac-038765e6fc83.md")
$ws.Hyperlinks.Add($ws.Range("C8"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c90dde53bac128dcff7171856fc5b578629fa35e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/978a0a30-878a-4dbc-a5ac-038765e6fc83.c933c29170ae9c133805761f650ea2014570545e.de-de.xlf", "", "", "978a0a30-878a-4dbc-a5ac-038765e6fc83.c933c29170ae9c133805761f650ea2014570545e.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A9"), "https://github.com/OpenLocalizationTest/oltest/blob/bdc60e2bda5902d2ba3712a83dc583abec29b99c/e2e/ef12ab35-0a62-4357-af99-1bc873e0fcc3.md", "", "", "ef12ab35-0a62-4357-af99-1bc873e0fcc3.md")
$ws.Hyperlinks.Add($ws.Range("C9"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a1f4c6e8b23d059172a6e4c8d0b3f7a9510cde6b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/ef12ab35-0a62-4357-af99-1bc873e0fcc3.b0a908d97a0c162733732f1a2a4d7c80cbd57c92.de-de.xlf", "", "", "ef12ab35-0a62-4357-af99-1bc873e0fcc3.b0a908d97a0c162733732f1a2a4d7c80cbd57c92.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A10"), "https://github.com/OpenLocalizationTest/oltest/blob/bdc60e2bda5902d2ba3712a83dc583abec29b99c/.localization-config", "", "", ".localization-config")

Write-Output "done"
